# Insert a new weekly price record as row 148 on the active sheet.
# This pushes the existing rows 148:216 down to 149:217 (dimension
# grows from A1:T216 to A1:T217) and populates the newly-inserted
# row 148 with the new observation. The "dimension"/category columns
# (A, B, C, E-L) carry the same values as the record that used to sit
# in row 148 (now shifted to row 149); only the date (D), the
# volume/price columns (M-P), the unit (Q), the origin (R), the
# per-kg price (S) and the kg/unit (T) differ for this new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 148:216 down by one row.
$ws.Rows("148:148").Insert()

# Populate the newly inserted row 148 with the new record.
$ws.Range("A148").Value = 6
$ws.Range("B148").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C148").Value = "Metropolitana"
$ws.Range("D148").Value = 44518
$ws.Range("E148").Value = 13
$ws.Range("F148").Value = "Fruta"
$ws.Range("G148").Value = 100101
$ws.Range("H148").Value = "Berries"
$ws.Range("I148").Value = 100101001
$ws.Range("J148").Value = "Arándano (blue)"
$ws.Range("K148").Value = "Sin especificar"
$ws.Range("L148").Value = "Primera"
$ws.Range("M148").Value = 2750
$ws.Range("N148").Value = 5000
$ws.Range("O148").Value = 5000
$ws.Range("P148").Value = 5000
$ws.Range("Q148").Value = "$/bandeja 2 kilos"
$ws.Range("R148").Value = "Región de O'Higgins"
$ws.Range("S148").Value = 2500
$ws.Range("T148").Value = 2
